# Auto-generated edit script: refreshes market-price-derived columns
# (H:N -> currentAveragePrice*, LevePrice*, LeveProfit*) for the rows
# touched by this scheduled data refresh, per sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 383.66666
$arr[0,1] = 75.5
$arr[0,2] = 1000
$arr[0,3] = 226.5
$arr[0,4] = 3000
$arr[0,5] = -114.5
$arr[0,6] = -3224
$ws.Range("H6:N6").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2900
$arr[0,1] = 2900
$arr[0,2] = 0
$arr[0,3] = 2900
$arr[0,4] = 0
$arr[0,5] = -2731
$arr[0,6] = $null
$ws.Range("H13:N13").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2874.75
$arr[0,1] = 2874.75
$arr[0,2] = 0
$arr[0,3] = 2874.75
$arr[0,4] = 0
$arr[0,5] = -2644.75
$arr[0,6] = $null
$ws.Range("H20:N20").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2874.75
$arr[0,1] = 2874.75
$arr[0,2] = 0
$arr[0,3] = 2874.75
$arr[0,4] = 0
$arr[0,5] = -2495.75
$arr[0,6] = $null
$ws.Range("H35:N35").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 13290118
$arr[0,1] = 325
$arr[0,2] = 14653173
$arr[0,3] = 975
$arr[0,4] = 43959519
$arr[0,5] = 133
$arr[0,6] = -43961735
$ws.Range("H112:N112").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 849.8
$arr[0,1] = 0
$arr[0,2] = 849.8
$arr[0,3] = 0
$arr[0,4] = 7648.2
$arr[0,5] = $null
$arr[0,6] = -12568.2
$ws.Range("H125:N125").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1255.4
$arr[0,1] = 708.0833
$arr[0,2] = 2076.375
$arr[0,3] = 2124.2499
$arr[0,4] = 6229.125
$arr[0,5] = 2875.7501
$arr[0,6] = -16229.125
$ws.Range("H129:N129").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2125.1392
$arr[0,1] = 1215.1892
$arr[0,2] = 2926.762
$arr[0,3] = 3645.5676
$arr[0,4] = 8780.286
$arr[0,5] = 1494.4324
$arr[0,6] = -19060.286
$ws.Range("H138:N138").Value = $arr

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 7590
$arr[0,1] = 1966.6666
$arr[0,2] = 10000
$arr[0,3] = 1966.6666
$arr[0,4] = 10000
$arr[0,5] = -1851.6666
$arr[0,6] = -10230
$ws.Range("H3:N3").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8000
$arr[0,1] = 8000
$arr[0,2] = 0
$arr[0,3] = 8000
$arr[0,4] = 0
$arr[0,5] = -7856
$arr[0,6] = $null
$ws.Range("H8:N8").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = 0
$ws.Range("H11:N11").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6524.149
$arr[0,1] = 5665.3257
$arr[0,2] = 21811.2
$arr[0,3] = 5665.3257
$arr[0,4] = 21811.2
$arr[0,5] = -5378.3257
$arr[0,6] = -22385.2
$ws.Range("H32:N32").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3268.375
$arr[0,1] = 1735.2858
$arr[0,2] = 14000
$arr[0,3] = 1735.2858
$arr[0,4] = 14000
$arr[0,5] = -1321.2858
$arr[0,6] = -14828
$ws.Range("H41:N41").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1823.25
$arr[0,1] = 1386.4762
$arr[0,2] = 2657.0908
$arr[0,3] = 1386.4762
$arr[0,4] = 2657.0908
$arr[0,5] = -512.4762000000001
$arr[0,6] = -4405.0908
$ws.Range("H74:N74").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1823.25
$arr[0,1] = 1386.4762
$arr[0,2] = 2657.0908
$arr[0,3] = 6932.381
$arr[0,4] = 13285.454
$arr[0,5] = -2564.381
$arr[0,6] = -22021.454
$ws.Range("H77:N77").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2225165.8
$arr[0,1] = 2062.7576
$arr[0,2] = 8338698.5
$arr[0,3] = 6188.2728
$arr[0,4] = 25016095.5
$arr[0,5] = -3658.2728
$arr[0,6] = -25021155.5
$ws.Range("H132:N132").Value = $arr

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 322.9
$arr[0,1] = 295
$arr[0,2] = 434.5
$arr[0,3] = 295
$arr[0,4] = 434.5
$arr[0,5] = -122
$arr[0,6] = -780.5
$ws.Range("H22:N22").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = 0
$ws.Range("H57:N57").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 40486
$arr[0,1] = 0
$arr[0,2] = 40486
$arr[0,3] = 0
$arr[0,4] = 40486
$arr[0,5] = $null
$arr[0,6] = -50606
$ws.Range("H133:N133").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = 0
$ws.Range("H136:N136").Value = $arr

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 54287
$arr[0,1] = 20000
$arr[0,2] = 100003
$arr[0,3] = 20000
$arr[0,4] = 100003
$arr[0,5] = -19887
$arr[0,6] = -100229
$ws.Range("H3:N3").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = 0
$ws.Range("H11:N11").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = 0
$ws.Range("H47:N47").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1609.4286
$arr[0,1] = 973.875
$arr[0,2] = 3643.2
$arr[0,3] = 2921.625
$arr[0,4] = 10929.6
$arr[0,5] = -391.625
$arr[0,6] = -15989.6
$ws.Range("H132:N132").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 192439.97
$arr[0,1] = 3988.8857
$arr[0,2] = 558872.6
$arr[0,3] = 11966.6571
$arr[0,4] = 1676617.8
$arr[0,5] = -9431.6571
$arr[0,6] = -1681687.8
$ws.Range("H134:N134").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 33160
$arr[0,1] = 0
$arr[0,2] = 33160
$arr[0,3] = 0
$arr[0,4] = 33160
$arr[0,5] = $null
$arr[0,6] = -43300
$ws.Range("H135:N135").Value = $arr

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3930.6
$arr[0,1] = 5984.778
$arr[0,2] = 1755.5883
$arr[0,3] = 17954.334
$arr[0,4] = 5266.7649
$arr[0,5] = -17842.334
$arr[0,6] = -5490.7649
$ws.Range("H5:N5").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 62500200
$arr[0,1] = 398
$arr[0,2] = 125000000
$arr[0,3] = 1194
$arr[0,4] = 375000000
$arr[0,5] = -970
$arr[0,6] = -375000448
$ws.Range("H9:N9").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 7751.4
$arr[0,1] = 784.75
$arr[0,2] = 15713.286
$arr[0,3] = 7062.75
$arr[0,4] = 141419.574
$arr[0,5] = -4612.75
$arr[0,6] = -146319.574
$ws.Range("H122:N122").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 25420.566
$arr[0,1] = 88274.164
$arr[0,2] = 7024.39
$arr[0,3] = 264822.492
$arr[0,4] = 21073.17
$arr[0,5] = -259762.492
$arr[0,6] = -31193.17
$ws.Range("H133:N133").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3930.6
$arr[0,1] = 5984.778
$arr[0,2] = 1755.5883
$arr[0,3] = 53863.002
$arr[0,4] = 15800.2947
$arr[0,5] = -51328.002
$arr[0,6] = -20870.2947
$ws.Range("H135:N135").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 19264.111
$arr[0,1] = 7052.2383
$arr[0,2] = 62005.668
$arr[0,3] = 21156.7149
$arr[0,4] = 186017.004
$arr[0,5] = -16056.7149
$arr[0,6] = -196217.004
$ws.Range("H137:N137").Value = $arr

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 180000
$arr[0,1] = 180000
$arr[0,2] = 0
$arr[0,3] = 180000
$arr[0,4] = 0
$arr[0,5] = -179830
$arr[0,6] = $null
$ws.Range("H9:N9").Value = $arr

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2650.9
$arr[0,1] = 2434
$arr[0,2] = 2976.25
$arr[0,3] = 2434
$arr[0,4] = 2976.25
$arr[0,5] = -2322
$arr[0,6] = -3200.25
$ws.Range("H7:N7").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 11900
$arr[0,1] = 0
$arr[0,2] = 11900
$arr[0,3] = 0
$arr[0,4] = 11900
$arr[0,5] = $null
$arr[0,6] = -12820
$ws.Range("H39:N39").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 90912216
$arr[0,1] = 111114376
$arr[0,2] = 2500
$arr[0,3] = 111114376
$arr[0,4] = 2500
$arr[0,5] = -111114240
$arr[0,6] = -2772
$ws.Range("H40:N40").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6501967
$arr[0,1] = 6501967
$arr[0,2] = 0
$arr[0,3] = 19505901
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = -19503451
$ws.Range("H122:N122").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2650.9
$arr[0,1] = 2434
$arr[0,2] = 2976.25
$arr[0,3] = 7302
$arr[0,4] = 8928.75
$arr[0,5] = -4832
$arr[0,6] = -13868.75
$ws.Range("H126:N126").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 12352598
$arr[0,1] = 15158668
$arr[0,2] = 5887.8
$arr[0,3] = 45476004
$arr[0,4] = 17663.4
$arr[0,5] = -45473474
$arr[0,6] = -22723.4
$ws.Range("H132:N132").Value = $arr

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 166671170
$arr[0,1] = 250002000
$arr[0,2] = 9503.5
$arr[0,3] = 250002000
$arr[0,4] = 9503.5
$arr[0,5] = -250001858
$arr[0,6] = -9787.5
$ws.Range("H12:N12").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 39809.668
$arr[0,1] = 0
$arr[0,2] = 39809.668
$arr[0,3] = 0
$arr[0,4] = 39809.668
$arr[0,5] = $null
$arr[0,6] = -40271.668
$ws.Range("H46:N46").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1840.1578
$arr[0,1] = 1280.9333
$arr[0,2] = 3937.25
$arr[0,3] = 3842.7999
$arr[0,4] = 11811.75
$arr[0,5] = -1312.7999
$arr[0,6] = -16871.75
$ws.Range("H132:N132").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 39809.668
$arr[0,1] = 0
$arr[0,2] = 39809.668
$arr[0,3] = 0
$arr[0,4] = 119429.004
$arr[0,5] = $null
$arr[0,6] = -124499.004
$ws.Range("H134:N134").Value = $arr
